# Reorder (reverse) the comma-separated "Recorded By" names in column G
# so that the first listed recorder becomes the last, e.g.
#   "backup@backdoor.com, System" -> "System, backup@backdoor.com"
# Applies to every used row on the active worksheet where column G contains
# a comma-separated list (i.e. more than one token).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $reversed = $parts[($parts.Length - 1)..0]
        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
